$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-7 down to 4-8
$ws.Rows("3:3").Insert()

# Seed the new row by duplicating row 4 (its formatting/text-storage),
# then overwrite every cell with the new match data below.
$ws.Rows("4:4").Copy()
$ws.Rows("3:3").PasteSpecial()

# Populate the newly inserted row 3 with the new match data (Sport Recife vs Ceara)
# NOTE: column B (Date) is left untouched - the row copy above already seeded it
# with the correct literal text "07/10/2024" (same value needed here), and
# re-assigning that string via .Value would make Excel re-parse it as a date.
$ws.Cells.Item(3, 1).Value = "dhuWUY6g"
$ws.Cells.Item(3, 3).Value = "21:00"
$ws.Cells.Item(3, 4).Value = "BRAZIL - SERIE B"
$ws.Cells.Item(3, 5).Value = "Sport Recife"
$ws.Cells.Item(3, 6).Value = "Ceara"
$ws.Cells.Item(3, 7).Value = 1.91
$ws.Cells.Item(3, 8).Value = 3.2
$ws.Cells.Item(3, 9).Value = 4.33
$ws.Cells.Item(3, 10).Value = 2.63
$ws.Cells.Item(3, 11).Value = 2.05
$ws.Cells.Item(3, 12).Value = 4.75
$ws.Cells.Item(3, 13).Value = 1.08
$ws.Cells.Item(3, 14).Value = 8
$ws.Cells.Item(3, 15).Value = 1.36
$ws.Cells.Item(3, 16).Value = 3
$ws.Cells.Item(3, 17).Value = 2.25
$ws.Cells.Item(3, 18).Value = 1.62
$ws.Cells.Item(3, 19).Value = 1.5
$ws.Cells.Item(3, 20).Value = 2.5
$ws.Cells.Item(3, 21).Value = 2
$ws.Cells.Item(3, 22).Value = 1.73
$ws.Cells.Item(3, 23).Value = 6
$ws.Cells.Item(3, 24).Value = 8.5
$ws.Cells.Item(3, 25).Value = 9
$ws.Cells.Item(3, 26).Value = 15
$ws.Cells.Item(3, 27).Value = 17
$ws.Cells.Item(3, 28).Value = 34
$ws.Cells.Item(3, 29).Value = 8
$ws.Cells.Item(3, 30).Value = 6.5
$ws.Cells.Item(3, 31).Value = 17
$ws.Cells.Item(3, 32).Value = 51
$ws.Cells.Item(3, 33).Value = 10
$ws.Cells.Item(3, 34).Value = 21
$ws.Cells.Item(3, 35).Value = 15
$ws.Cells.Item(3, 36).Value = 41
$ws.Cells.Item(3, 37).Value = 41
$ws.Cells.Item(3, 38).Value = 41
$ws.Cells.Item(3, 39).Value = 351
$ws.Cells.Item(3, 40).Value = 3.75
$ws.Cells.Item(3, 41).Value = 11
$ws.Cells.Item(3, 42).Value = 23
$ws.Cells.Item(3, 43).Value = 41
$ws.Cells.Item(3, 44).Value = 67
$ws.Cells.Item(3, 45).Value = 201
$ws.Cells.Item(3, 46).Value = 2.5
$ws.Cells.Item(3, 47).Value = 8.5
$ws.Cells.Item(3, 48).Value = 67
$ws.Cells.Item(3, 49).Value = 6
$ws.Cells.Item(3, 50).Value = 23
$ws.Cells.Item(3, 51).Value = 34
$ws.Cells.Item(3, 52).Value = 81
$ws.Cells.Item(3, 53).Value = 126
$ws.Cells.Item(3, 54).Value = 301
$ws.Cells.Item(3, 55).Value = 126
$ws.Cells.Item(3, 56).Value = 126

# Minor odds corrections on row 4 (previously row 3) that changed beyond the shift
$ws.Cells.Item(4, 7).Value = 2.2
$ws.Cells.Item(4, 9).Value = 3.2
$ws.Cells.Item(4, 12).Value = 3.5
$ws.Cells.Item(4, 33).Value = 12
$ws.Cells.Item(4, 34).Value = 17
$ws.Cells.Item(4, 40).Value = 4.5
$ws.Cells.Item(4, 41).Value = 12
$ws.Cells.Item(4, 49).Value = 5
$ws.Cells.Item(4, 51).Value = 21

# Minor odds corrections on row 5 (previously row 4) that changed beyond the shift
$ws.Cells.Item(5, 13).Value = 1.07
$ws.Cells.Item(5, 14).Value = 9
$ws.Cells.Item(5, 17).Value = 2.25
$ws.Cells.Item(5, 18).Value = 1.62
